$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
Write-Host $ws.Name
